$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CambioTasa")

# Update the data row (row 2) with the new transaction values
$ws.Range("A2").Value = "btorres"
$ws.Range("C2").Value = "'1008535937"
$ws.Range("D2").Value = "'8"
$ws.Range("E2").Value = "ACHACALTANAS1"
$ws.Range("G2").Value = "AAACT232011SD7TZT"
$ws.Range("H2").Value = "20 jul. 2023, 09:42:54"

# Update the last selected cell, matching the saved view state
$ws.Range("E10").Select()
